$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-08 Thursday" "2024-02-09 Friday"
Replace-Text "30×16=480" "25×64=1600"
Replace-Text "16×66=1056" "50×21=1050"
Replace-Text "96×77=7392" "80×47=3760"
Replace-Text "77×27=2079" "39×80=3120"
Replace-Text "50×73=3650" "39×91=3549"
Replace-Text "81×11=891" "94×81=7614"
Replace-Text "72×82=5904" "63×50=3150"
Replace-Text "87×66=5742" "39×71=2769"
Replace-Text "86×25=2150" "86×24=2064"
Replace-Text "28×23=644" "39×89=3471"
Replace-Text "11×66=726" "47×86=4042"
Replace-Text "66×98=6468" "90×24=2160"
Replace-Text "54×65=3510" "51×54=2754"
Replace-Text "75×47=3525" "20×77=1540"
Replace-Text "66×94=6204" "99×42=4158"
Replace-Text "45×61=2745" "91×12=1092"
Replace-Text "32×79=2528" "57×98=5586"
Replace-Text "15×61=915" "47×93=4371"
Replace-Text "25×29=725" "48×17=816"
Replace-Text "35×26=910" "47×90=4230"
Replace-Text "19×16=304" "91×12=1092"
Replace-Text "52×58=3016" "30×86=2580"
Replace-Text "59×38=2242" "85×46=3910"
Replace-Text "36×27=972" "94×99=9306"
Replace-Text "39×46=1794" "60×97=5820"
